$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original plain-text representation
# (values like "236.10" or "42.79" would otherwise be auto-converted to
# floating point numbers by Excel, losing trailing zeros / exact formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.076.95"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").Value = "1.855.55"
$ws.Range("E3").Value = "  +2.72%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "236.10"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("E6").Value = "  +2.09%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "42.79"
$ws.Range("E8").Value = "  +9.57%  "

$ws.Range("E9").Value = "  +2.47%  "

$ws.Range("E10").Value = "  +2.33%  "

$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").Value = "2.125.28"
$ws.Range("E12").Value = "  +2.72%  "

$ws.Range("E13").Value = "  +2.36%  "

$ws.Range("D14").Value = "1.855.39"
$ws.Range("E14").Value = "  +2.64%  "

$ws.Range("D15").Value = "0.680"
$ws.Range("E15").Value = "  +2.90%  "

$ws.Range("D16").Value = "4.69"
$ws.Range("E16").Value = "  +2.79%  "

$ws.Range("D17").Value = "35.021.23"
$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("D18").Value = "70.33"
$ws.Range("E18").Value = "  +1.82%  "

$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  +2.00%  "

$ws.Range("D20").Value = "241.06"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value = "  +3.17%  "

$ws.Range("E22").Value = "  +2.12%  "

$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  +2.63%  "

$ws.Range("D25").Value = "171.78"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  +24.92%  "

$ws.Range("D27").Value = "7.91"
$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("D28").Value = "17.67"
$ws.Range("E28").Value = "  +2.79%  "

$ws.Range("D29").Value = "0.124"
$ws.Range("E29").Value = "  +2.66%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  +2.93%  "

$ws.Range("E32").Value = "  -1.24%  "

$ws.Range("D33").Value = "3.99"
$ws.Range("E33").Value = "  +2.74%  "

$ws.Range("D34").Value = "2.03"
$ws.Range("E34").Value = "  +13.73%  "

$ws.Range("E35").Value = "  +22.85%  "

$ws.Range("D36").Value = "0.780"
$ws.Range("E36").Value = "  +12.64%  "

$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("E38").Value = "  +13.12%  "

$ws.Range("D39").Value = "91.66"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("E40").Value = "  +6.00%  "

$ws.Range("D41").Value = "1.350.22"
$ws.Range("E41").Value = "  +2.25%  "

$ws.Range("E42").Value = "  +5.78%  "

$ws.Range("E43").Value = "  +6.02%  "

$ws.Range("D44").Value = "12.77"
$ws.Range("E44").Value = "  +61.13%  "

$ws.Range("E45").Value = "  -2.27%  "

$ws.Range("D46").Value = "2.76"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("D47").Value = "0.0548"
$ws.Range("E47").Value = "  +6.94%  "

$ws.Range("D48").Value = "6.43"
$ws.Range("E48").Value = "  +3.50%  "

$ws.Range("D49").Value = "2.039.47"
$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("D50").Value = "0.0681"
$ws.Range("E50").Value = "  +2.74%  "

$ws.Range("D51").Value = "3.41"
$ws.Range("E51").Value = "  +15.21%  "
